# "Avance muestreo de datos"
# Add three new rows (CC, Pasaporte, TI) of sample data to the
# "TipoIdentificacion" sheet, and make that sheet the active/selected one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TipoIdentificacion")

# --- Add the new data rows (6, 7, 8 -> rows 7, 8, 9) -------------------
$newRows = @(
    @{ Row = 7; Num = 6; Nombre = "CC" },
    @{ Row = 8; Num = 7; Nombre = "Pasaporte" },
    @{ Row = 9; Num = 8; Nombre = "TI" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Num
    $ws.Cells.Item($r, 2).Value = $item.Nombre
    $ws.Cells.Item($r, 3).Formula = "=B" + $r
}

# Copy the formatting (styles/borders/fill) of the last existing data row
# down onto the three new rows, matching the look of the table above.
$ws.Range("A6:C6").Copy() | Out-Null
$ws.Range("A7:C9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Make "TipoIdentificacion" the active sheet/tab ---------------------
$ws.Activate() | Out-Null
$ws.Range("F8:F9").Select() | Out-Null
